$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "n" counts on row 1 (columns B:E)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update the CON row (row 2) values for columns B:E
$ws.Range("B2").Value = 2.9921861609128104
$ws.Range("C2").Value = 0.6791466556503124
$ws.Range("D2").Value = 2.8514699578225438
$ws.Range("E2").Value = 0.38689077983612485

# Update the STR row (row 3) values for columns B:E
$ws.Range("B3").Value = 2.4520921249786043
$ws.Range("C3").Value = 0.91898135611724596
$ws.Range("D3").Value = 2.7423587664283593
$ws.Range("E3").Value = 0.67348917426278099

# Update the active selection to match the edited range
$ws.Range("B1:E3").Select()
